# Update "想去人数" (interested-count) values on the "展览" and "全部类型"
# sheets to reflect freshly scraped counts (gh-pages data refresh @456a3b4).

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> { row -> new F-column value }
$updates = @{
    "展览" = @{
        3  = 21528
        8  = 8010
        9  = 560
        12 = 321
        18 = 233
        19 = 1363
        20 = 553
        21 = 86
        24 = 90
        27 = 1204
        28 = 62
        33 = 152
        34 = 5123
        38 = 60
        39 = 13239
        40 = 1372
        41 = 147
        42 = 59
        44 = 328
        45 = 453
        46 = 4072
        47 = 16
        49 = 104
    }
    "全部类型" = @{
        3  = 21528
        6  = 8010
        7  = 560
        10 = 321
        15 = 233
        16 = 1363
        17 = 553
        18 = 86
        21 = 90
        24 = 1204
        25 = 62
        32 = 152
        34 = 5123
        38 = 60
        39 = 13239
        40 = 1372
        41 = 147
        42 = 59
        44 = 328
        45 = 453
        46 = 4072
        47 = 16
        49 = 104
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item($row, 6).Value = $rows[$row]
    }
}
